# Update header labels (D1, E1) to include units
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "E (MPa)"
$ws.Range("E1").Value = "k (N/mm)"

# Re-populate the data rows (the raw test-data rows were shuffled/re-sorted
# and the tissue/team/E/k values for each row were updated accordingly).

$ws.Range("B2").Value = "LCL"
$ws.Range("C2").Value = "'02"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 4.833454170321591
$ws.Range("E2").Value = 1.325462690005715

$ws.Range("B3").Value = "MCL"
$ws.Range("C3").Value = "'08"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 8.612915666753445
$ws.Range("E3").Value = 2.351042034748962

$ws.Range("B4").Value = "LCL"
$ws.Range("C4").Value = "'08"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = 10.22701617502049
$ws.Range("E4").Value = 1.363602156669399

$ws.Range("B5").Value = "MCL"
$ws.Range("C5").Value = "'05"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = 20.04660225334274
$ws.Range("E5").Value = 4.747879481054859

$ws.Range("B6").Value = "MCL"
$ws.Range("C6").Value = "'04"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 17.60619541581704
$ws.Range("E6").Value = 5.711278025130895

$ws.Range("B7").Value = "LCL"
$ws.Range("C7").Value = "'10"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = 19.89518258065792
$ws.Range("E7").Value = 4.50456964090368

$ws.Range("B8").Value = "MCL"
$ws.Range("C8").Value = "'09"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 21.75123279164934
$ws.Range("E8").Value = 4.857775323468353

$ws.Range("B9").Value = "MCL"
$ws.Range("C9").Value = "'10"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = 9.433230463407448
$ws.Range("E9").Value = 3.074534373258724

$ws.Range("B10").Value = "LCL"
$ws.Range("C10").Value = "'04"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = 21.07680203522363
$ws.Range("E10").Value = 2.965673914690758

$ws.Range("B11").Value = "LCL"
$ws.Range("C11").Value = "'03"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = 13.32303589345438
$ws.Range("E11").Value = 4.946177075444941

$ws.Range("B12").Value = "MCL"
$ws.Range("C12").Value = "'07"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = 23.92605810181052
$ws.Range("E12").Value = 5.638021612109806

$ws.Range("B13").Value = "LCL"
$ws.Range("C13").Value = "'01"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = 12.11483330268918
$ws.Range("E13").Value = 4.330408499684641

$ws.Range("B14").Value = "MCL"
$ws.Range("C14").Value = "'02"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = 3.510014164356793
$ws.Range("E14").Value = 1.864159961679736

$ws.Range("B15").Value = "MCL"
$ws.Range("C15").Value = "'06"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = 32.40808937901795
$ws.Range("E15").Value = 6.328105873481926

$ws.Range("B16").Value = "LCL"
$ws.Range("C16").Value = "'05"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 31.7608406328724
$ws.Range("E16").Value = 4.466368213997681

$ws.Range("B17").Value = "MCL"
$ws.Range("C17").Value = "'01"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = 2.378275167186548
$ws.Range("E17").Value = 0.9409697400607646

$ws.Range("B18").Value = "LCL"
$ws.Range("C18").Value = "'09"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = 10.8491727760429
$ws.Range("E18").Value = 2.608770553017949

$ws.Range("B19").Value = "MCL"
$ws.Range("C19").Value = "'03"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = 12.41061955319874
$ws.Range("E19").Value = 8.211409923554783

$ws.Range("B20").Value = "LCL"
$ws.Range("C20").Value = "'06"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = 2.264103576822557
$ws.Range("E20").Value = 0.3069970951623806

$ws.Range("B21").Value = "LCL"
$ws.Range("C21").Value = "'07"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = 9.168396981614865
$ws.Range("E21").Value = 1.384132189482502

